$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: change B27 from time-of-day decimal 0.75 to text "18.00-20.00",
# and add D27/E27 notes, plus G27 hours value.
$ws.Range("B27").Value = "18.00-20.00"
$ws.Range("D27").Value = "Vähän maanantai ja heikot unet alla pitkän päivän iltana, mutta ehk siitä jotain jäi käteen. Erityisesti kiinnostaa katsoa ymmärsinkö kaiken oleellisen mallintaakseni seuraavaa demoa."
$ws.Range("E27").Value = "Tässähän alkaa tätä teknistä velkaakin jo olla, ehkä sitten kurssin loppupuolella katsellaan jos jää aikaa."
$ws.Range("G27").Value = 2

# Match D/E column wrap-text formatting already used elsewhere (e.g. D3/E3)
$ws.Range("D27").WrapText = $true
$ws.Range("E27").WrapText = $true

$ws.Rows.Item(27).RowHeight = 58

# Row 28: new entry for "12 marras"
$ws.Range("A28").Value = "12 marras"
$ws.Range("B28").Value = "14.00-16.00"
$ws.Range("C28").Value = "Kovien kappaleiden demo"
$ws.Range("G28").Value = 2

$ws.Range("B28").NumberFormat = $ws.Range("B26").NumberFormat
$ws.Range("B28").WrapText = $true
$ws.Range("C28").WrapText = $true

$ws.Range("H28").Select()
